# Update crypto price/volume data per the Sat Jun 10 06:13:19 UTC 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.659.70'
$ws.Range("E2").Value = '  -3.21%  '
# Row 3
$ws.Range("D3").Value = '1.764.06'
$ws.Range("E3").Value = '  -4.14%  '
# Row 4
$ws.Range("D4").Value = '''1.007'
$ws.Range("E4").Value = '  +0.57%  '
# Row 5
$ws.Range("D5").Value = '''237.64'
$ws.Range("E5").Value = '  -8.69%  '
# Row 6
$ws.Range("D6").Value = '''1.005'
$ws.Range("E6").Value = '  +0.40%  '
# Row 7
$ws.Range("D7").Value = '''0.4886'
$ws.Range("E7").Value = '  -6.96%  '
# Row 8
$ws.Range("D8").Value = '''42.30'
$ws.Range("E8").Value = '  -5.67%  '
# Row 9
$ws.Range("D9").Value = '''0.2486'
$ws.Range("E9").Value = '  -21.66%  '
# Row 10
$ws.Range("D10").Value = '''0.06068'
$ws.Range("E10").Value = '  -10.57%  '
# Row 11
$ws.Range("D11").Value = '1.784.67'
$ws.Range("E11").Value = '  -2.58%  '
# Row 12
$ws.Range("D12").Value = '''0.06594'
$ws.Range("E12").Value = '  -15.03%  '
# Row 13
$ws.Range("D13").Value = '''14.43'
$ws.Range("E13").Value = '  -22.74%  '
# Row 14
$ws.Range("D14").Value = '''0.6041'
$ws.Range("E14").Value = '  -22.71%  '
# Row 15
$ws.Range("D15").Value = '''78.00'
$ws.Range("E15").Value = '  -11.22%  '
# Row 16
$ws.Range("D16").Value = '''4.325'
$ws.Range("E16").Value = '  -13.68%  '
# Row 17
$ws.Range("D17").Value = '''1.007'
$ws.Range("E17").Value = '  +0.60%  '
# Row 18
$ws.Range("D18").Value = '''1.006'
$ws.Range("E18").Value = '  +0.51%  '
# Row 19
$ws.Range("D19").Value = '25.695.66'
$ws.Range("E19").Value = '  -3.14%  '
# Row 20
$ws.Range("D20").Value = '''11.05'
$ws.Range("E20").Value = '  -20.16%  '
# Row 21
$ws.Range("D21").Value = '2.008.68'
$ws.Range("E21").Value = '  -3.09%  '
# Row 22
$ws.Range("D22").Value = '''0.000006265'
$ws.Range("E22").Value = '  -21.24%  '
# Row 23
$ws.Range("D23").Value = '''3.870'
$ws.Range("E23").Value = '  -16.09%  '
# Row 24
$ws.Range("D24").Value = '''5.115'
$ws.Range("E24").Value = '  -14.39%  '
# Row 25
$ws.Range("D25").Value = '''8.023'
$ws.Range("E25").Value = '  -13.96%  '
# Row 26
$ws.Range("D26").Value = '''131.90'
$ws.Range("E26").Value = '  -7.47%  '
# Row 27
$ws.Range("D27").Value = '''1.875'
$ws.Range("E27").Value = '  -15.04%  '
# Row 28
$ws.Range("D28").Value = '''14.42'
$ws.Range("E28").Value = '  -14.70%  '
# Row 29
$ws.Range("D29").Value = '''1.382'
$ws.Range("E29").Value = '  -17.61%  '
# Row 30
$ws.Range("D30").Value = '''99.42'
$ws.Range("E30").Value = '  -10.97%  '
# Row 31
$ws.Range("D31").Value = '''0.08198'
$ws.Range("E31").Value = '  -5.93%  '
# Row 32
$ws.Range("D32").Value = '''3.577'
$ws.Range("E32").Value = '  -14.14%  '
# Row 33
$ws.Range("E33").Value = '  +0.60%  '
# Row 34
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''3.162'
$ws.Range("E34").Value = '  -22.30%  '
# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.04271'
$ws.Range("E35").Value = '  -12.61%  '
# Row 36
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '''2.627'
$ws.Range("E36").Value = '  -8.11%  '
# Row 37
$ws.Range("D37").Value = '''1.020'
$ws.Range("E37").Value = '  -9.88%  '
# Row 38
$ws.Range("D38").Value = '''0.6140'
$ws.Range("E38").Value = '  -15.39%  '
# Row 39
$ws.Range("D39").Value = '''2.657'
$ws.Range("E39").Value = '  -14.11%  '
# Row 40
$ws.Range("D40").Value = '''2.091'
$ws.Range("E40").Value = '  -6.19%  '
# Row 41
$ws.Range("D41").Value = '''1.005'
$ws.Range("E41").Value = '  +0.37%  '
# Row 42
$ws.Range("D42").Value = '''101.38'
$ws.Range("E42").Value = '  -7.87%  '
# Row 43
$ws.Range("D43").Value = '''0.01434'
$ws.Range("E43").Value = '  -17.87%  '
# Row 44
$ws.Range("D44").Value = '''0.7849'
$ws.Range("E44").Value = '  -12.36%  '
# Row 45
$ws.Range("D45").Value = '''0.3824'
# Row 46
$ws.Range("D46").Value = '''5.146'
$ws.Range("E46").Value = '  -12.99%  '
# Row 47
$ws.Range("D47").Value = '''6.098'
$ws.Range("E47").Value = '  -20.07%  '
# Row 48
$ws.Range("D48").Value = '''0.05179'
$ws.Range("E48").Value = '  -11.32%  '
# Row 49
$ws.Range("D49").Value = '''52.14'
$ws.Range("E49").Value = '  -12.51%  '
# Row 50
$ws.Range("D50").Value = '''1.004'
$ws.Range("E50").Value = '  -0.04%  '
# Row 51
$ws.Range("E51").Value = '  +0.21%  '
